$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Devices")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# --- Row 8: Built-in Loop-A (entered first, matches shared-string order) ---
$ws.Range("F8").Value = "Built-in Loop-A"

# --- Row 7 (header): new columns F (Loop) and G (Column Number) ---
$ws.Range("F7").Value = "Loop"
$ws.Range("G7").Value = "Column Number"

# --- Row 9: Built-in Loop-B ---
$ws.Range("F9").Value = "Built-in Loop-B"

# --- Row 8: Column Number value ---
$ws.Range("G8").Value = 1

# --- Copy header/data formats onto the new cells ---
$ws.Range("E7").Copy()
$ws.Range("F7:G7").PasteSpecial(-4122)

$ws.Range("E8").Copy()
$ws.Range("F8:G8").PasteSpecial(-4122)

$ws.Range("E9").Copy()
$ws.Range("F9").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# --- Column widths for new columns F and G (target OOXML width 13.21875 / 15) ---
$ws.Columns.Item(6).ColumnWidth = 12.333333333333334
$ws.Columns.Item(7).ColumnWidth = 14.166666666666666

# --- Update selection to reflect new active cell ---
$ws.Range("F11").Select()
